$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C11"  = -13.0642
    "B12"  = 4.727699999999997
    "D14"  = -8.559799999999996
    "D19"  = -8.2766
    "C23"  = -12.1276
    "D24"  = -8.37619999999999
    "B27"  = 5.624200000000002
    "C28"  = -13.5844
    "B32"  = 6.445599999999995
    "C32"  = -11.78060000000001
    "C34"  = -11.77320000000001
    "B36"  = 9.378900000000002
    "B38"  = 6.172700000000003
    "D38"  = -7.909500000000006
    "D41"  = -8.425799999999999
    "C42"  = -12.5256
    "B46"  = 6.089200000000002
    "C49"  = -13.43069999999999
    "D52"  = -7.805000000000003
    "B54"  = 5.259900000000001
    "C54"  = -14.14639999999999
    "B55"  = 5.589699999999996
    "B56"  = 4.811899999999997
    "B67"  = 5.494099999999995
    "B69"  = 5.290599999999995
    "B72"  = 5.225500000000004
    "D72"  = -7.243799999999998
    "C78"  = -11.78770000000001
    "D78"  = -7.844000000000003
    "C80"  = -13.2923
    "B83"  = 5.379199999999996
    "D83"  = -9.212799999999993
    "D85"  = -8.680400000000002
    "B86"  = 5.418899999999999
    "D86"  = -7.870699999999998
    "D90"  = -7.257099999999993
    "B91"  = 4.721499999999995
    "B93"  = 5.527400000000003
    "D96"  = -8.652399999999993
    "C97"  = -11.3669
    "B99"  = 5.977899999999998
    "C99"  = -12.43160000000001
    "C101" = -13.1941
    "D103" = -8.014500000000002
    "B104" = 9.607999999999999
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$wb.Save()
